# Actualización desde MV -datos-
# Append 4 new daily rows (02-11-2021 .. 05-11-2021) to the "Spot 2021 - Diaria" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates = @("02-11-2021", "03-11-2021", "04-11-2021", "05-11-2021")
$data = @(
    @(6.02, 5.07, 3.98, 3.69),
    @(6.12, 5.1,  4.01, 3.78),
    @(6.27, 5.17, 4.07, 3.88),
    @(6.23, 5.09, 4.02, 3.85)
)

$startRow = 213
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $startRow + $i

    # Writing the date string directly through Range.Value would be
    # re-interpreted as a real date (it's ambiguous as mm-dd-yyyy), changing
    # the cell's type/style. Instead enter it as a literal-text formula
    # (no number-like auto-conversion, no style change) and then convert
    # that formula to a plain value via copy / paste-special-values so it
    # lands as a normal shared-string cell, matching the rest of column A.
    $ws.Cells.Item($row, 1).Formula = "=""" + $dates[$i] + """"
    $ws.Range("A" + $row).Copy()
    $ws.Range("A" + $row).PasteSpecial(-4163)

    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
    $ws.Cells.Item($row, 4).Value = $data[$i][2]
    $ws.Cells.Item($row, 5).Value = $data[$i][3]
}
